# Update "Overall" sheet (row 2: Simulation 1 summary stats)
$wb = $excel.ActiveWorkbook
$wsOverall = $wb.Worksheets.Item("Overall")

$wsOverall.Range("B2").Value = 67
$wsOverall.Range("C2").Value = 2
$wsOverall.Range("D2").Value = 1.7747296719856605
$wsOverall.Range("E2").Value = 0.33739693645722207
$wsOverall.Range("F2").Value = 2.4003921568627455
$wsOverall.Range("G2").Value = 52
$wsOverall.Range("H2").Value = 29
$wsOverall.Range("I2").Value = 81
$wsOverall.Range("J2").Value = 575
$wsOverall.Range("K2").Value = 75

# Update "Zones" sheet (rows 2-14: per-zone stats)
$wsZones = $wb.Worksheets.Item("Zones")

# Row 2 (Zone 1)
$wsZones.Range("B2").Value = 9
$wsZones.Range("C2").Value = 1
$wsZones.Range("D2").Value = 0.51999999999999991
$wsZones.Range("E2").Value = 0.54074074074074074
$wsZones.Range("F2").Value = 0.33333333333333326

# Row 3 (Zone 2) - E3 becomes empty
$wsZones.Range("B3").Value = 5
$wsZones.Range("D3").Value = 1.6261904761904762
$wsZones.Range("E3").ClearContents()
$wsZones.Range("F3").Value = 1.6261904761904762

# Row 4 (Zone 3)
$wsZones.Range("B4").Value = 2
$wsZones.Range("D4").Value = 2.6777777777777776
$wsZones.Range("E4").Value = 0.28333333333333321
$wsZones.Range("F4").Value = 3.3619047619047615

# Row 5 (Zone 4)
$wsZones.Range("B5").Value = 5
$wsZones.Range("D5").Value = 1.6041666666666661
$wsZones.Range("E5").Value = 0.21666666666666648
$wsZones.Range("F5").Value = 2.4366666666666661

# Row 6 (Zone 5) - B6, C6 unchanged
$wsZones.Range("D6").Value = 2.0944444444444446
$wsZones.Range("E6").Value = 0.26111111111111107
$wsZones.Range("F6").Value = 3.0111111111111115

# Row 7 (Zone 6)
$wsZones.Range("B7").Value = 1
$wsZones.Range("D7").Value = 1.7708333333333333
$wsZones.Range("E7").Value = 0.21666666666666679
$wsZones.Range("F7").Value = 1.9928571428571427

# Row 8 (Zone 7)
$wsZones.Range("B8").Value = 9
$wsZones.Range("C8").Value = 0
$wsZones.Range("D8").Value = 2.3576923076923082
$wsZones.Range("E8").Value = 0.26666666666666672
$wsZones.Range("F8").Value = 3.2870370370370376

# Row 9 (Zone 8)
$wsZones.Range("B9").Value = 7
$wsZones.Range("D9").Value = 1.3222222222222226
$wsZones.Range("E9").Value = 0.22333333333333333
$wsZones.Range("F9").Value = 2.1071428571428577

# Row 10 (Zone 9) - E10 newly added
$wsZones.Range("B10").Value = 1
$wsZones.Range("C10").Value = 1
$wsZones.Range("D10").Value = 1.2194444444444443
$wsZones.Range("E10").Value = 0.66666666666666674
$wsZones.Range("F10").Value = 1.4958333333333331

# Row 11 (Zone 10)
$wsZones.Range("B11").Value = 1
$wsZones.Range("D11").Value = 2.9999999999999991
$wsZones.Range("E11").Value = 0.1416666666666665
$wsZones.Range("F11").Value = 3.5196969696969687

# Row 12 (Zone 11)
$wsZones.Range("B12").Value = 17
$wsZones.Range("D12").Value = 0.95784313725490189
$wsZones.Range("E12").Value = 0.23999999999999982
$wsZones.Range("F12").Value = 1.2569444444444444

# Row 13 (Zone 12)
$wsZones.Range("B13").Value = 2
$wsZones.Range("C13").Value = 0
$wsZones.Range("D13").Value = 1.3729166666666668
$wsZones.Range("E13").Value = 0.41666666666666696
$wsZones.Range("F13").Value = 1.5095238095238095

# Row 14 (Zone 13) - E14 becomes empty
$wsZones.Range("B14").Value = 0
$wsZones.Range("D14").Value = 4.95
$wsZones.Range("E14").ClearContents()
$wsZones.Range("F14").Value = 4.95
